$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 72 as a duplicate of the current (last) data row 71,
# preserving that row's original values (this becomes the new row 72).
$ws.Rows.Item(71).Copy()
$ws.Rows.Item(72).Insert()

# Update row 71 with this week's new data (newer date and updated prices).
$ws.Range("D71").Value = 44509
$ws.Range("K71").Value = 7500
$ws.Range("L71").Value = 8000
$ws.Range("M71").Value = 7750
$ws.Range("P71").Value = 155
